$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New LeetCode rows 114-118. Copy the format (styles, column widths, number
# formats) from the last existing data row (113) down into each new row so
# the new cells pick up the same style indexes (centered / wrap-text /
# date) as the rest of the table, then overwrite the copied values.

# Row 114: Rearranging Fruits (hard)
$ws.Range("A113:I113").Copy($ws.Range("A114:I114"))
$ws.Rows.Item(114).RowHeight = 34
$ws.Cells.Item(114, 1).Value2 = 2561
$ws.Cells.Item(114, 2).Value2 = "Rearranging Fruits"
$ws.Cells.Item(114, 3).Value2 = "#greedy #hash-table #sorting "
$ws.Cells.Item(114, 4).Value2 = "hard"
$ws.Cells.Item(114, 5).Value2 = 0
$ws.Cells.Item(114, 6).Value2 = 1
$ws.Cells.Item(114, 7).Value2 = 20
$ws.Cells.Item(114, 8).Value2 = 45871
$ws.Cells.Item(114, 9).Value2 = 45871

# Row 115: One Edit Distance (medium)
$ws.Range("A113:I113").Copy($ws.Range("A115:I115"))
$ws.Rows.Item(115).RowHeight = 34
$ws.Cells.Item(115, 1).Value2 = 161
$ws.Cells.Item(115, 2).Value2 = "One Edit Distance"
$ws.Cells.Item(115, 3).Value2 = "#string #two-pointers #array "
$ws.Cells.Item(115, 4).Value2 = "medium"
$ws.Cells.Item(115, 5).Value2 = 0
$ws.Cells.Item(115, 6).Value2 = 1
$ws.Cells.Item(115, 7).Value2 = 13
$ws.Cells.Item(115, 8).Value2 = 45872
$ws.Cells.Item(115, 9).Value2 = 45872

# Row 116: Fruit Into Baskets (medium)
$ws.Range("A113:I113").Copy($ws.Range("A116:I116"))
$ws.Rows.Item(116).RowHeight = 34
$ws.Cells.Item(116, 1).Value2 = 904
$ws.Cells.Item(116, 2).Value2 = "Fruit Into Baskets"
$ws.Cells.Item(116, 3).Value2 = "#array #sliding-window "
$ws.Cells.Item(116, 4).Value2 = "medium"
$ws.Cells.Item(116, 5).Value2 = 0
$ws.Cells.Item(116, 6).Value2 = 1
$ws.Cells.Item(116, 7).Value2 = 20
$ws.Cells.Item(116, 8).Value2 = 45873
$ws.Cells.Item(116, 9).Value2 = 45873

# Row 117: Fruits Into Baskets II (easy)
$ws.Range("A113:I113").Copy($ws.Range("A117:I117"))
$ws.Rows.Item(117).RowHeight = 17
$ws.Cells.Item(117, 1).Value2 = 3477
$ws.Cells.Item(117, 2).Value2 = "Fruits Into Baskets II"
$ws.Cells.Item(117, 3).Value2 = "#array"
$ws.Cells.Item(117, 4).Value2 = "easy"
$ws.Cells.Item(117, 5).Value2 = 1
$ws.Cells.Item(117, 6).Value2 = 0
$ws.Cells.Item(117, 7).Value2 = 5
$ws.Cells.Item(117, 8).Value2 = 45874
$ws.Cells.Item(117, 9).Value2 = 45874

# Row 118: Maximum Fruits Harvested After at Most K Steps (hard)
$ws.Range("A113:I113").Copy($ws.Range("A118:I118"))
$ws.Rows.Item(118).RowHeight = 51
$ws.Cells.Item(118, 1).Value2 = 2106
$ws.Cells.Item(118, 2).Value2 = "Maximum Fruits Harvested After at Most K Steps"
$ws.Cells.Item(118, 3).Value2 = "#array #two-pointers #sliding-window "
$ws.Cells.Item(118, 4).Value2 = "hard"
$ws.Cells.Item(118, 5).Value2 = 0
$ws.Cells.Item(118, 6).Value2 = 1
$ws.Cells.Item(118, 7).Value2 = 36
$ws.Cells.Item(118, 8).Value2 = 45874
$ws.Cells.Item(118, 9).Value2 = 45874

# Match the author's final cursor position/selection (H118:I118).
$ws.Range("H118:I118").Select()
